# Update table 1.16.B from "October" YTD figures to "November" YTD figures
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table_1_16_B")

# --- Title / header text updates -------------------------------------------------
$ws.Range("A2").Value = "by State, by Sector, Year-to-Date through November 2016 and 2015 (Thousand Megawatthours)"

# Column header labels (shared across B/E/G/I/K and C/F/H/J/L via row 6)
$ws.Range("B6").Value = "November 2016 YTD"
$ws.Range("E6").Value = "November 2016 YTD"
$ws.Range("G6").Value = "November 2016 YTD"
$ws.Range("I6").Value = "November 2016 YTD"
$ws.Range("K6").Value = "November 2016 YTD"

$ws.Range("C6").Value = "November 2015 YTD"
$ws.Range("F6").Value = "November 2015 YTD"
$ws.Range("H6").Value = "November 2015 YTD"
$ws.Range("J6").Value = "November 2015 YTD"
$ws.Range("L6").Value = "November 2015 YTD"

# --- Data updates ------------------------------------------------------------------
# Row 52: Mountain
$ws.Range("B52").Value = 4064
$ws.Range("C52").Value = 3252
$ws.Range("D52").Value = 0.25
$ws.Range("E52").Value = 234
$ws.Range("F52").Value = 237
$ws.Range("G52").Value = 3830
$ws.Range("H52").Value = 3015

# Row 55: Idaho
$ws.Range("B55").Value = 88
$ws.Range("C55").Value = 68
$ws.Range("D55").Value = 0.296
$ws.Range("G55").Value = 88
$ws.Range("H55").Value = 68

# Row 57: Nevada
$ws.Range("B57").Value = 3482
$ws.Range("C57").Value = 2784
$ws.Range("D57").Value = 0.251
$ws.Range("G57").Value = 3482
$ws.Range("H57").Value = 2784

# Row 58: New Mexico
$ws.Range("B58").Value = 14
$ws.Range("C58").Value = 9
$ws.Range("D58").Value = 0.656
$ws.Range("G58").Value = 14
$ws.Range("H58").Value = 9

# Row 59: Utah
$ws.Range("B59").Value = 480
$ws.Range("C59").Value = 391
$ws.Range("D59").Value = 0.228
$ws.Range("E59").Value = 234
$ws.Range("F59").Value = 237
$ws.Range("G59").Value = 245
$ws.Range("H59").Value = 154

# Row 61: Pacific Contiguous
$ws.Range("B61").Value = 11501
$ws.Range("C61").Value = 11079
$ws.Range("D61").Value = 0.038
$ws.Range("E61").Value = 763
$ws.Range("F61").Value = 754
$ws.Range("G61").Value = 10738
$ws.Range("H61").Value = 10325

# Row 62: California
$ws.Range("B62").Value = 11327
$ws.Range("C62").Value = 10920
$ws.Range("D62").Value = 0.037
$ws.Range("E62").Value = 747
$ws.Range("F62").Value = 753
$ws.Range("G62").Value = 10580
$ws.Range("H62").Value = 10167

# Row 63: Oregon
$ws.Range("B63").Value = 174
$ws.Range("C63").Value = 159
$ws.Range("D63").Value = 0.095
$ws.Range("F63").Value = 1
$ws.Range("G63").Value = 158
$ws.Range("H63").Value = 158

# Row 65: Pacific Noncontiguous
$ws.Range("B65").Value = 231
$ws.Range("C65").Value = 210
$ws.Range("D65").Value = 0.102
$ws.Range("G65").Value = 231
$ws.Range("H65").Value = 210

# Row 67: Hawaii
$ws.Range("B67").Value = 231
$ws.Range("C67").Value = 210
$ws.Range("D67").Value = 0.102
$ws.Range("G67").Value = 231
$ws.Range("H67").Value = 210

# Row 68: U.S. Total
$ws.Range("B68").Value = 15797
$ws.Range("C68").Value = 14541
$ws.Range("D68").Value = 0.086
$ws.Range("E68").Value = 997
$ws.Range("F68").Value = 991
$ws.Range("G68").Value = 14800
$ws.Range("H68").Value = 13549
